$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.583.10'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.859.66'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.68%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.99'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4678'
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3892'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.29'
$ws.Range("E9").Value = '  -5.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07974'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9994'
$ws.Range("E11").Value = '  -3.03%  '
$ws.Range("E12").Value = '  -2.37%  '
$ws.Range("D13").Value = '1.873.16'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.969'
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.239'
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.014'
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.99'
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06725'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001041'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.95'
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("D22").Value = '27.575.98'
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.431'
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.84'
$ws.Range("E24").Value = '  -1.72%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '2.084.39'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.77'
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.72'
$ws.Range("E28").Value = '  -2.30%  '
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.388'
$ws.Range("E30").Value = '  -3.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.16'
$ws.Range("E31").Value = '  -0.57%  '
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09450'
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.278'
$ws.Range("E35").Value = '  -1.48%  '
$ws.Range("E36").Value = '  -8.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06031'
$ws.Range("E37").Value = '  -1.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02219'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.190'
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.221'
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.011'
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5907'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1873'
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.250'
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5613'
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.13'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.278'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06761'
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.38'
$ws.Range("E51").Value = '  -1.71%  '
